# Applies the "Add files via upload" edit:
#  - Slide 2 (CONTENTS): drop the "IMAGE RECOGNITION" line, renumber the
#    remaining lines, and nudge the bullet-list placeholder's position.
#  - Slide 7 ("...RESULT:" title): renumber 6 -> 5.
#  - Slide 8 ("...ADVANTAGES AND DISADVANTAGES" title): renumber 7 -> 6.
#  - Slide 9 ("...CONCLUSION" title): renumber 8 -> 7.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - CONTENTS list
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)          # "Subtitle 2" placeholder holding the list

# --- reposition the placeholder -----------------------------------------
# Target offsets (EMU): x=4797393, y=1518377. EMU/point = 12700. The COM
# Left/Top properties round-trip through a lower-precision float, so nudge
# iteratively until the EMU we get back on save matches exactly.
$emuPerPt = 12700.0
$targetXEmu = 4797393
$targetYEmu = 1518377
$xPt = $targetXEmu / $emuPerPt
$yPt = $targetYEmu / $emuPerPt
for ($i = 0; $i -lt 6; $i++) {
    $sh2.Left = $xPt
    $sh2.Top = $yPt
    $curXEmu = [math]::Round($sh2.Left * $emuPerPt)
    $curYEmu = [math]::Round($sh2.Top * $emuPerPt)
    if ($curXEmu -eq $targetXEmu -and $curYEmu -eq $targetYEmu) { break }
    $xPt = $xPt + (($targetXEmu - $curXEmu) / $emuPerPt)
    $yPt = $yPt + (($targetYEmu - $curYEmu) / $emuPerPt)
}

# --- update the text content ---------------------------------------------
$tr2 = $sh2.TextFrame.TextRange

# Remove the old paragraph 8 ("8.CONCLUSION") outright - the surrounding
# paragraphs shift up and the trailing empty paragraph is preserved as-is.
$tr2.Paragraphs(8, 1).Delete()

# Rename the remaining paragraphs (from the end up, so earlier character
# offsets stay valid while each edit is applied).
$tr2.Characters(100, 31).Text = "7.CONCLUSION"                      # was 7.ADVANTAGES AND DISADVANTAGES
$tr2.Characters(91, 9).Text  = "6.ADVANTAGES AND DISADVANTAGES"     # was 6.RESULT
$tr2.Characters(71, 20).Text = "5.RESULT"                           # was 5.IMAGE RECOGNITION

# ---------------------------------------------------------------------
# Slide 7 - "6.RESULT:" -> "5.RESULT:"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(1)          # Title
$tr7 = $sh7.TextFrame.TextRange
$tr7.Characters(1, 1).Text = "5"

# ---------------------------------------------------------------------
# Slide 8 - "7. ADVANTAGES AND DISADVANTAGES" -> "6. ..."
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(1)          # Title
$tr8 = $sh8.TextFrame.TextRange
$tr8.Characters(1, 1).Text = "6"

# ---------------------------------------------------------------------
# Slide 9 - "8.CONCLUSION" -> "7.CONCLUSION"
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)          # Title
$tr9 = $sh9.TextFrame.TextRange
$tr9.Text = "7.CONCLUSION"
